# Oklahoma overview workbook restructuring:
#  - Reorders the data columns so "Share of 990 filers with government
#    grants at risk" becomes the first data column (right after any
#    Geography/Size/Subsector label column).
#  - Renames the "Operating surplus with/without government grants (%)"
#    headers to "Size of operating surplus with/without government
#    grants".
#  - Renames several row labels (Congressional districts, Size buckets,
#    Subsectors) and reorders the Size-bucket and Subsector rows.
#
# Because values like "906", "72.63%", "$1,608,045,127" must stay as
# literal text (not get reinterpreted as numbers/percentages by Excel),
# every destination cell gets NumberFormat "@" (Text) applied before the
# new value is written.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

function Get-CellText($ws, $row, $col) {
    return $ws.Cells.Item($row, $col).Value2
}

# ---------------------------------------------------------------------
# Sheet "Overall": 5 data columns, A:E, header row 1 + single data row 2.
# Old order: A=Number, B=Total$, C=OpSurplusWith%, D=OpSurplusWithout%,
#            E=ShareAtRisk
# New order: A=ShareAtRisk, B=Number, C=Total$,
#            D=SizeOpSurplusWith, E=SizeOpSurplusWithout
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")

$o_h_A = Get-CellText $wsOverall 1 1
$o_h_B = Get-CellText $wsOverall 1 2
$o_h_C = Get-CellText $wsOverall 1 3
$o_h_D = Get-CellText $wsOverall 1 4
$o_h_E = Get-CellText $wsOverall 1 5

$o_d_A = Get-CellText $wsOverall 2 1
$o_d_B = Get-CellText $wsOverall 2 2
$o_d_C = Get-CellText $wsOverall 2 3
$o_d_D = Get-CellText $wsOverall 2 4
$o_d_E = Get-CellText $wsOverall 2 5

$wsOverall.Range("A1:E2").NumberFormat = "@"

Set-TextCell $wsOverall 1 1 $o_h_E
Set-TextCell $wsOverall 1 2 $o_h_A
Set-TextCell $wsOverall 1 3 $o_h_B
Set-TextCell $wsOverall 1 4 "Size of operating surplus with government grants"
Set-TextCell $wsOverall 1 5 "Size of operating surplus without government grants"

Set-TextCell $wsOverall 2 1 $o_d_E
Set-TextCell $wsOverall 2 2 $o_d_A
Set-TextCell $wsOverall 2 3 $o_d_B
Set-TextCell $wsOverall 2 4 $o_d_C
Set-TextCell $wsOverall 2 5 $o_d_D

# ---------------------------------------------------------------------
# Generic helper for the six-column sheets (County, Congressional
# District, Size, Subsector): column A is the label, B:F are data.
# Old order: B=Number, C=Total$, D=OpSurplusWith%, E=OpSurplusWithout%,
#            F=ShareAtRisk
# New order: B=ShareAtRisk, C=Number, D=Total$, E=SizeOpSurplusWith,
#            F=SizeOpSurplusWithout
# ---------------------------------------------------------------------
function Set-SixColHeader($ws) {
    $h_A = Get-CellText $ws 1 1
    $h_B = Get-CellText $ws 1 2
    $h_C = Get-CellText $ws 1 3

    $ws.Range("A1:F1").NumberFormat = "@"

    Set-TextCell $ws 1 1 $h_A
    Set-TextCell $ws 1 2 "Share of 990 filers with government grants at risk"
    Set-TextCell $ws 1 3 $h_B
    Set-TextCell $ws 1 4 $h_C
    Set-TextCell $ws 1 5 "Size of operating surplus with government grants"
    Set-TextCell $ws 1 6 "Size of operating surplus without government grants"
}

# Reads the old B:F values of $srcRow and writes the shifted values into
# $dstRow (A is left untouched here - callers set the label separately).
function Copy-ShiftedRow($ws, $srcRow, $dstRow) {
    $b = Get-CellText $ws $srcRow 2
    $c = Get-CellText $ws $srcRow 3
    $d = Get-CellText $ws $srcRow 4
    $e = Get-CellText $ws $srcRow 5
    $f = Get-CellText $ws $srcRow 6

    $ws.Range("B$dstRow`:F$dstRow").NumberFormat = "@"

    Set-TextCell $ws $dstRow 2 $f
    Set-TextCell $ws $dstRow 3 $b
    Set-TextCell $ws $dstRow 4 $c
    Set-TextCell $ws $dstRow 5 $d
    Set-TextCell $ws $dstRow 6 $e
}

# ---------------------------------------------------------------------
# Sheet "County": label column unchanged, 68 data rows (2-69), no
# reordering or relabeling - only the column shift applies. Read every
# row's B:F values up-front (rows are written in place so src==dst, but
# we still read first for consistency/safety).
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
Set-SixColHeader $wsCounty

$countyOld = @{}
for ($r = 2; $r -le 69; $r++) {
    $countyOld[$r] = @(
        (Get-CellText $wsCounty $r 2),
        (Get-CellText $wsCounty $r 3),
        (Get-CellText $wsCounty $r 4),
        (Get-CellText $wsCounty $r 5),
        (Get-CellText $wsCounty $r 6)
    )
}
for ($r = 2; $r -le 69; $r++) {
    $vals = $countyOld[$r]
    $b = $vals[0]; $c = $vals[1]; $d = $vals[2]; $e = $vals[3]; $f = $vals[4]
    $wsCounty.Range("B$r`:F$r").NumberFormat = "@"
    Set-TextCell $wsCounty $r 2 $f
    Set-TextCell $wsCounty $r 3 $b
    Set-TextCell $wsCounty $r 4 $c
    Set-TextCell $wsCounty $r 5 $d
    Set-TextCell $wsCounty $r 6 $e
}

# ---------------------------------------------------------------------
# Sheet "Congressional District": rows 2 (US) & 3 (Oklahoma) keep their
# position; rows 4-8 keep their position but the label changes from
# "Nth Congressional district" to "Congressional District N".
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
Set-SixColHeader $wsCd

# Snapshot B:F for all data rows before writing (src == dst here, but
# keep the read/write phases separate to mirror the other sheets).
$cdOld = @{}
for ($r = 2; $r -le 8; $r++) {
    $cdOld[$r] = @(
        (Get-CellText $wsCd $r 2),
        (Get-CellText $wsCd $r 3),
        (Get-CellText $wsCd $r 4),
        (Get-CellText $wsCd $r 5),
        (Get-CellText $wsCd $r 6)
    )
}
for ($r = 2; $r -le 8; $r++) {
    $vals = $cdOld[$r]
    $b = $vals[0]; $c = $vals[1]; $d = $vals[2]; $e = $vals[3]; $f = $vals[4]
    $wsCd.Range("B$r`:F$r").NumberFormat = "@"
    Set-TextCell $wsCd $r 2 $f
    Set-TextCell $wsCd $r 3 $b
    Set-TextCell $wsCd $r 4 $c
    Set-TextCell $wsCd $r 5 $d
    Set-TextCell $wsCd $r 6 $e
}

$cdLabels = @{
    4 = "Congressional District 1";
    5 = "Congressional District 2";
    6 = "Congressional District 3";
    7 = "Congressional District 4";
    8 = "Congressional District 5";
}
foreach ($r in $cdLabels.Keys) {
    $wsCd.Cells.Item($r, 1).NumberFormat = "@"
    $wsCd.Cells.Item($r, 1).Value = $cdLabels[$r]
}

# ---------------------------------------------------------------------
# Sheet "Size": rows reordered into ascending size order, with the
# following new-row <- old-row data mapping (label rewritten too):
#   2 <- 2  "Between $100K and $499K"
#   3 <- 4  "Between $1M and $4.99M"
#   4 <- 5  "Between $500K and $999K"
#   5 <- 6  "Between $5M and $9.99M"
#   6 <- 3  "Greater than $10M"
#   7 <- 7  "Less than $100K"            (unchanged position)
#   8 <- 8  "Total"                       (unchanged position)
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
Set-SixColHeader $wsSize

# Snapshot every old row (A label + B:F) before any write, since the
# rows move around and some are overwritten before being read otherwise.
$sizeOld = @{}
for ($r = 2; $r -le 8; $r++) {
    $sizeOld[$r] = @(
        (Get-CellText $wsSize $r 1),
        (Get-CellText $wsSize $r 2),
        (Get-CellText $wsSize $r 3),
        (Get-CellText $wsSize $r 4),
        (Get-CellText $wsSize $r 5),
        (Get-CellText $wsSize $r 6)
    )
}

$sizeRowMap = @{ 2 = 2; 3 = 4; 4 = 5; 5 = 6; 6 = 3; 7 = 7; 8 = 8 }
$sizeLabels = @{
    2 = "Between `$100K and `$499K";
    3 = "Between `$1M and `$4.99M";
    4 = "Between `$500K and `$999K";
    5 = "Between `$5M and `$9.99M";
    6 = "Greater than `$10M";
    7 = "Less than `$100K";
    8 = "Total";
}

foreach ($dstRow in 2..8) {
    $srcRow = $sizeRowMap[$dstRow]
    $vals = $sizeOld[$srcRow]
    $b = $vals[1]; $c = $vals[2]; $d = $vals[3]; $e = $vals[4]; $f = $vals[5]

    $wsSize.Range("A$dstRow`:F$dstRow").NumberFormat = "@"
    Set-TextCell $wsSize $dstRow 1 $sizeLabels[$dstRow]
    Set-TextCell $wsSize $dstRow 2 $f
    Set-TextCell $wsSize $dstRow 3 $b
    Set-TextCell $wsSize $dstRow 4 $c
    Set-TextCell $wsSize $dstRow 5 $d
    Set-TextCell $wsSize $dstRow 6 $e
}

# ---------------------------------------------------------------------
# Sheet "Subsector": rows 2-10 and 13 keep their position (label text is
# tweaked for several of them); rows 11 and 12 swap positions (and the
# "Universities"/"Unclassified" labels swap accordingly).
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
Set-SixColHeader $wsSub

$subOld = @{}
for ($r = 2; $r -le 13; $r++) {
    $subOld[$r] = @(
        (Get-CellText $wsSub $r 1),
        (Get-CellText $wsSub $r 2),
        (Get-CellText $wsSub $r 3),
        (Get-CellText $wsSub $r 4),
        (Get-CellText $wsSub $r 5),
        (Get-CellText $wsSub $r 6)
    )
}

$subRowMap = @{ 2 = 2; 3 = 3; 4 = 4; 5 = 5; 6 = 6; 7 = 7; 8 = 8; 9 = 9; 10 = 10; 11 = 12; 12 = 11; 13 = 13 }
$subLabels = @{
    2 = "Arts, Culture, and Humanities";
    3 = "Education (Excluding Universities)";
    4 = "Environment and Animals";
    5 = "Health (Excluding Hospitals)";
    6 = "Hospitals";
    7 = "Human Services";
    8 = "International, Foreign Affairs";
    9 = "Public, Societal Benefit";
    10 = "Religion Related";
    11 = "Unclassified";
    12 = "Universities";
    13 = "Total";
}

foreach ($dstRow in 2..13) {
    $srcRow = $subRowMap[$dstRow]
    $vals = $subOld[$srcRow]
    $b = $vals[1]; $c = $vals[2]; $d = $vals[3]; $e = $vals[4]; $f = $vals[5]

    $wsSub.Range("A$dstRow`:F$dstRow").NumberFormat = "@"
    Set-TextCell $wsSub $dstRow 1 $subLabels[$dstRow]
    Set-TextCell $wsSub $dstRow 2 $f
    Set-TextCell $wsSub $dstRow 3 $b
    Set-TextCell $wsSub $dstRow 4 $c
    Set-TextCell $wsSub $dstRow 5 $d
    Set-TextCell $wsSub $dstRow 6 $e
}
